# Refresh the cryptos list data (Price / Volume(1h) columns, plus a few
# rows whose coin name+link moved to a different rank position) to match
# the latest GitHub Actions scrape.
#
# Column D ("Price") cells are stored as literal text in this sheet (even
# when the text happens to look like a plain number, e.g. "4.34"), so a
# leading apostrophe is used for those to force Excel to keep them as text
# instead of silently converting them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20's price contains a unicode subscript-three character (U+2083).
$sub3 = [char]0x2083


# Row 2
$ws.Range("D2").Value = "34.628.54"
$ws.Range("E2").Value = "  +1.75%  "

# Row 3
$ws.Range("D3").Value = "1.836.52"
$ws.Range("E3").Value = "  +3.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.59%  "

# Row 5
$ws.Range("D5").Formula = "'225.92"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("D6").Formula = "'0.555"
$ws.Range("E6").Value = "  +1.22%  "

# Row 7
$ws.Range("E7").Value = "  -0.63%  "

# Row 8
$ws.Range("D8").Formula = "'32.75"
$ws.Range("E8").Value = "  +5.20%  "

# Row 9
$ws.Range("E9").Value = "  +4.33%  "

# Row 10
$ws.Range("D10").Formula = "'0.0711"
$ws.Range("E10").Value = "  +8.31%  "

# Row 11
$ws.Range("D11").Formula = "'0.0931"
$ws.Range("E11").Value = "  +0.03%  "

# Row 12
$ws.Range("D12").Value = "2.093.16"
$ws.Range("E12").Value = "  +2.67%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.833.85"
$ws.Range("E13").Value = "  +3.06%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Formula = "'11.17"
$ws.Range("E14").Value = "  +0.22%  "

# Row 15
$ws.Range("E15").Value = "  +4.29%  "

# Row 16
$ws.Range("D16").Value = "34.639.44"
$ws.Range("E16").Value = "  +1.76%  "

# Row 17
$ws.Range("D17").Formula = "'4.34"
$ws.Range("E17").Value = "  +3.17%  "

# Row 18
$ws.Range("D18").Formula = "'69.71"
$ws.Range("E18").Value = "  +1.14%  "

# Row 19
$ws.Range("D19").Formula = "'254.01"
$ws.Range("E19").Value = "  +0.71%  "

# Row 20
$ws.Range("D20").Value = "0.0{0}0802" -f $sub3
$ws.Range("E20").Value = "  +8.64%  "

# Row 21
$ws.Range("D21").Formula = "'11.36"
$ws.Range("E21").Value = "  +9.34%  "

# Row 22
$ws.Range("E22").Value = "  -0.35%  "

# Row 23
$ws.Range("D23").Formula = "'4.30"
$ws.Range("E23").Value = "  +2.46%  "

# Row 24
$ws.Range("E24").Value = "  +1.47%  "

# Row 25
$ws.Range("D25").Formula = "'161.89"
$ws.Range("E25").Value = "  +3.75%  "

# Row 26
$ws.Range("D26").Formula = "'16.76"
$ws.Range("E26").Value = "  +2.09%  "

# Row 27
$ws.Range("D27").Formula = "'7.22"
$ws.Range("E27").Value = "  +3.28%  "

# Row 28
$ws.Range("E28").Value = "  +1.02%  "

# Row 29
$ws.Range("D29").Formula = "'0.996"
$ws.Range("E29").Value = "  -0.72%  "

# Row 30
$ws.Range("D30").Formula = "'0.0531"
$ws.Range("E30").Value = "  +3.06%  "

# Row 31
$ws.Range("E31").Value = "  +1.75%  "

# Row 32
$ws.Range("D32").Formula = "'1.21"
$ws.Range("E32").Value = "  +0.66%  "

# Row 33
$ws.Range("D33").Formula = "'499.78"
$ws.Range("E33").Value = "  +864.50%  "

# Row 34
$ws.Range("E34").Value = "  +2.42%  "

# Row 35
$ws.Range("D35").Formula = "'1.94"
$ws.Range("E35").Value = "  +6.27%  "

# Row 36
$ws.Range("D36").Value = "1.440.29"
$ws.Range("E36").Value = "  -0.39%  "

# Row 37
$ws.Range("D37").Formula = "'0.656"
$ws.Range("E37").Value = "  +4.65%  "

# Row 38
$ws.Range("E38").Value = "  +1.47%  "

# Row 39
$ws.Range("D39").Formula = "'0.0193"
$ws.Range("E39").Value = "  +3.32%  "

# Row 40
$ws.Range("D40").Formula = "'0.977"
$ws.Range("E40").Value = "  +9.69%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Formula = "'2.86"
$ws.Range("E41").Value = "  +0.66%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Formula = "'83.22"
$ws.Range("E42").Value = "  +0.67%  "

# Row 43
$ws.Range("D43").Formula = "'2.37"
$ws.Range("E43").Value = "  +0.54%  "

# Row 44
$ws.Range("E44").Value = "  +5.16%  "

# Row 45
$ws.Range("D45").Formula = "'6.11"
$ws.Range("E45").Value = "  +5.49%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.989.87"
$ws.Range("E46").Value = "  +2.69%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Formula = "'12.51"
$ws.Range("E47").Value = "  +6.28%  "

# Row 48
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Formula = "'1.06"
$ws.Range("E48").Value = "  -0.04%  "

# Row 49
$ws.Range("D49").Formula = "'0.0493"
$ws.Range("E49").Value = "  -3.14%  "

# Row 50
$ws.Range("D50").Formula = "'106.56"
$ws.Range("E50").Value = "  +9.48%  "

# Row 51
$ws.Range("E51").Value = "  -0.22%  "
